# Updates cryptos list price (D) / volume-1h (E) columns in place,
# matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D prices that look like plain numbers ("251.46", "0.640", ...) are
# forced to text with a leading apostrophe (quote-prefix) so Excel keeps them
# as literal strings -- exactly like the source data (e.g. preserves "0.640"
# instead of normalising it to 0.64, and keeps multi-dot values like
# "42.929.16" intact, which can never parse as a number anyway).

$ws.Range("D2").Value = '42.929.16'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '2.288.75'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'" + '251.46'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").Value = "'" + '0.640'
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("D7").Value = "'" + '73.43'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'" + '0.650'
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("D10").Value = "'" + '39.03'
$ws.Range("E10").Value = '  -5.44%  '
$ws.Range("D11").Value = "'" + '0.0976'
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("D12").Value = "'" + '59.09'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").Value = '2.630.38'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").Value = "'" + '15.29'
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = "'" + '0.872'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").Value = '2.290.62'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").Value = '42.835.73'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("D22").Value = "'" + '72.51'
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").Value = "'" + '234.84'
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").Value = '  +6.29%  '
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = "'" + '2.42'
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E30").Value = '  -3.42%  '
$ws.Range("D31").Value = "'" + '167.18'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("D33").Value = "'" + '6.43'
$ws.Range("E33").Value = '  +4.81%  '
$ws.Range("E34").Value = '  -3.71%  '
$ws.Range("D35").Value = "'" + '0.0827'
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("D36").Value = "'" + '31.03'
$ws.Range("E36").Value = '  +7.56%  '
$ws.Range("E37").Value = '  +1.54%  '
$ws.Range("D38").Value = "'" + '4.57'
$ws.Range("E38").Value = '  +9.84%  '
$ws.Range("E39").Value = '  +1.72%  '
$ws.Range("E40").Value = '  -3.99%  '
$ws.Range("E41").Value = '  +12.55%  '
$ws.Range("E42").Value = '  +2.08%  '
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  +7.19%  '
$ws.Range("D45").Value = "'" + '9.15'
$ws.Range("E45").Value = '  +2.28%  '
$ws.Range("D46").Value = "'" + '61.88'
$ws.Range("E46").Value = '  -4.40%  '
$ws.Range("D47").Value = "'" + '4.88'
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("D49").Value = "'" + '102.28'
$ws.Range("E49").Value = '  +7.99%  '
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").Value = "'" + '1.17'
$ws.Range("E51").Value = '  -1.93%  '
